$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on price cells whose new values would
# otherwise be auto-parsed as numbers, so they stay text like the rest
# of the Price column (matches source data which is all inline text).
$textCells = @("D5","D6","D8","D10","D18","D20","D23","D24","D25","D31","D37","D38","D39","D41","D43","D46","D47","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '26.008.18'
$ws.Range("E2").Value = '  -0.04%  '
$ws.Range("D3").Value = '1.634.87'
$ws.Range("E3").Value = '  -0.50%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '214.07'
$ws.Range("E5").Value = '  -1.03%  '
$ws.Range("D6").Value = '0.504'
$ws.Range("E6").Value = '  -0.58%  '
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D8").Value = '0.252'
$ws.Range("E8").Value = '  -2.10%  '
$ws.Range("E9").Value = '  -2.36%  '
$ws.Range("D10").Value = '18.49'
$ws.Range("E10").Value = '  -5.91%  '
$ws.Range("E11").Value = '  -0.42%  '
$ws.Range("D12").Value = '1.861.86'
$ws.Range("E13").Value = '  -1.94%  '
$ws.Range("D14").Value = '1.636.49'
$ws.Range("E14").Value = '  +1.27%  '
$ws.Range("E15").Value = '  -2.65%  '
$ws.Range("D16").Value = '26.007.92'
$ws.Range("E16").Value = '  +0.19%  '
$ws.Range("D17").Value = '0.0₃0745'
$ws.Range("E17").Value = '  -2.81%  '
$ws.Range("D18").Value = '61.81'
$ws.Range("E19").Value = '  +0.10%  '
$ws.Range("D20").Value = '190.31'
$ws.Range("E20").Value = '  -1.44%  '
$ws.Range("E21").Value = '  -2.28%  '
$ws.Range("E22").Value = '  -3.51%  '
$ws.Range("D23").Value = '6.14'
$ws.Range("E23").Value = '  -2.02%  '
$ws.Range("D24").Value = '0.133'
$ws.Range("E24").Value = '  +0.48%  '
$ws.Range("D25").Value = '143.24'
$ws.Range("E25").Value = '  -0.93%  '
$ws.Range("E27").Value = '  -2.12%  '
$ws.Range("E28").Value = '  -2.74%  '
$ws.Range("E29").Value = '  -2.23%  '
$ws.Range("E30").Value = '  -1.42%  '
$ws.Range("D31").Value = '0.0483'
$ws.Range("E31").Value = '  -3.37%  '
$ws.Range("E32").Value = '  -2.64%  '
$ws.Range("E33").Value = '  -3.94%  '
$ws.Range("E34").Value = '  -1.56%  '
$ws.Range("E35").Value = '  -2.24%  '
$ws.Range("B36").Value = 'Maker'
$ws.Range("C36").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D36").Value = '1.135.27'
$ws.Range("E36").Value = '  +0.02%  '
$ws.Range("B37").Value = 'ARBITRUM'
$ws.Range("C37").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D37").Value = '0.870'
$ws.Range("E37").Value = '  -3.73%  '
$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D38").Value = '2.43'
$ws.Range("E38").Value = '  -1.51%  '
$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").Value = '0.526'
$ws.Range("E39").Value = '  -3.18%  '
$ws.Range("E40").Value = '  -1.47%  '
$ws.Range("D41").Value = '98.63'
$ws.Range("E41").Value = '  -1.14%  '
$ws.Range("E42").Value = '  -1.75%  '
$ws.Range("D43").Value = '5.27'
$ws.Range("E43").Value = '  -4.42%  '
$ws.Range("D44").Value = '1.772.09'
$ws.Range("E45").Value = '  -1.01%  '
$ws.Range("D46").Value = '55.26'
$ws.Range("E46").Value = '  -2.71%  '
$ws.Range("D47").Value = '0.0526'
$ws.Range("E47").Value = '  -0.65%  '
$ws.Range("E48").Value = '  +1.75%  '
$ws.Range("E49").Value = '  -0.40%  '
$ws.Range("D50").Value = '7.54'
$ws.Range("E50").Value = '  -3.19%  '
$ws.Range("D51").Value = '1.00'
$ws.Range("E51").Value = '  +0.04%  '
